# Updated cryptos list on Sat Jan  6 10:34:49 UTC 2024 with GitHub Actions
# Refresh the "Price" (D) and "Volume(1h)" (E) columns with the latest
# scraped values, and fix the ImmutableX/Dai row ordering (rows 25-26
# were swapped upstream).
#
# Note: several "Price" values are plain numeric-looking strings (e.g.
# "7.20", "0.522") that must stay as TEXT (matching the sheet's existing
# inline-string cells), so those are written with a leading apostrophe to
# force text entry and avoid Excel's automatic number coercion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.982.53"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.235.93"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'305.39"
$ws.Range("E5").Value = "  -4.13%  "
$ws.Range("D6").Value = "'95.47"
$ws.Range("E6").Value = "  -6.01%  "
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("E8").Value = "  +0.15%  "
$ws.Range("D9").Value = "'0.522"
$ws.Range("E9").Value = "  -5.61%  "
$ws.Range("D10").Value = "'34.95"
$ws.Range("E10").Value = "  -5.93%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  -3.59%  "
$ws.Range("D12").Value = "'7.20"
$ws.Range("E12").Value = "  -5.32%  "
$ws.Range("E13").Value = "  -2.85%  "
$ws.Range("D14").Value = "2.576.86"
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").Value = "2.241.31"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("D16").Value = "'0.824"
$ws.Range("E16").Value = "  -4.27%  "
$ws.Range("D17").Value = "'13.59"
$ws.Range("E17").Value = "  -6.65%  "
$ws.Range("D18").Value = "43.863.74"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "0.0₃0960"
$ws.Range("E19").Value = "  -2.79%  "
$ws.Range("D20").Value = "'12.19"
$ws.Range("E20").Value = "  -9.03%  "
$ws.Range("D21").Value = "'6.21"
$ws.Range("E21").Value = "  -4.24%  "
$ws.Range("D22").Value = "'64.89"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").Value = "'236.44"
$ws.Range("E23").Value = "  +0.46%  "
$ws.Range("E24").Value = "  -6.46%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.94"
$ws.Range("E25").Value = "  -6.73%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'9.92"
$ws.Range("E27").Value = "  -6.26%  "
$ws.Range("D28").Value = "'38.00"
$ws.Range("E28").Value = "  -1.58%  "
$ws.Range("D29").Value = "'2.15"
$ws.Range("E29").Value = "  -2.02%  "
$ws.Range("D30").Value = "'5.94"
$ws.Range("E30").Value = "  -4.17%  "
$ws.Range("D31").Value = "'19.88"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "'154.94"
$ws.Range("E32").Value = "  -4.45%  "
$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("E34").Value = "  +4.24%  "
$ws.Range("E35").Value = "  -3.41%  "
$ws.Range("E36").Value = "  -1.54%  "
$ws.Range("E37").Value = "  -5.81%  "
$ws.Range("E38").Value = "  -10.82%  "
$ws.Range("D39").Value = "'15.24"
$ws.Range("E39").Value = "  -8.51%  "
$ws.Range("D40").Value = "'3.35"
$ws.Range("E40").Value = "  -8.76%  "
$ws.Range("D41").Value = "'3.82"
$ws.Range("E41").Value = "  -8.93%  "
$ws.Range("D42").Value = "'0.0301"
$ws.Range("E42").Value = "  -4.93%  "
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.740.67"
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").Value = "'85.20"
$ws.Range("E45").Value = "  +4.22%  "
$ws.Range("D46").Value = "'0.188"
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("D47").Value = "'99.73"
$ws.Range("E47").Value = "  -4.68%  "
$ws.Range("D48").Value = "'4.91"
$ws.Range("E48").Value = "  -5.97%  "
$ws.Range("D49").Value = "'69.12"
$ws.Range("E49").Value = "  -7.34%  "
$ws.Range("D50").Value = "'8.07"
$ws.Range("E50").Value = "  -3.39%  "
$ws.Range("D51").Value = "'54.25"
$ws.Range("E51").Value = "  -6.84%  "
